$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Rows 293-304: years (column C) corrected from 6 to 3
$ws.Range("C293").Value = 3
$ws.Range("C294").Value = 3
$ws.Range("C295").Value = 3
$ws.Range("C296").Value = 3
$ws.Range("C297").Value = 3
$ws.Range("C298").Value = 3
$ws.Range("C299").Value = 3
$ws.Range("C300").Value = 3
$ws.Range("C301").Value = 3
$ws.Range("C302").Value = 3
$ws.Range("C303").Value = 3
$ws.Range("C304").Value = 3

# Rows 305-325: type (B), years (C) re-coded; se_higher/se_lower (E,G) switched to
# literal 1/SE-based precision values instead of the old D/100 and F/100 formulas;
# a few beta_lower (F) values corrected to point at the right preceding row
$ws.Range("C305").Value = 6
$ws.Range("E305").Value = 0.004

$ws.Range("B306").Value = 2
$ws.Range("C306").Value = 2
$ws.Range("E306").Value = 0.004
$ws.Range("G306").Value = 0.004

$ws.Range("B307").Value = 2
$ws.Range("C307").Value = 4
$ws.Range("E307").Value = 0.005
$ws.Range("F307").Value = 0.0968
$ws.Range("G307").Value = 0.004

$ws.Range("C308").Value = 6
$ws.Range("E308").Value = 0.004
$ws.Range("G308").Value = 0.005

$ws.Range("B309").Value = 2
$ws.Range("C309").Value = 2
$ws.Range("E309").Value = 0.004
$ws.Range("G309").Value = 0.004

$ws.Range("B310").Value = 2
$ws.Range("C310").Value = 4
$ws.Range("E310").Value = 0.005
$ws.Range("F310").Value = 0.10300000000000001
$ws.Range("G310").Value = 0.004

$ws.Range("C311").Value = 6
$ws.Range("E311").Value = 0.004
$ws.Range("G311").Value = 0.005

$ws.Range("B312").Value = 2
$ws.Range("C312").Value = 2
$ws.Range("E312").Value = 0.004
$ws.Range("G312").Value = 0.004

$ws.Range("B313").Value = 2
$ws.Range("C313").Value = 4
$ws.Range("E313").Value = 0.005
$ws.Range("F313").Value = 0.09939999999999999
$ws.Range("G313").Value = 0.004

$ws.Range("C314").Value = 6
$ws.Range("E314").Value = 0.01
$ws.Range("G314").Value = 0.005

$ws.Range("B315").Value = 2
$ws.Range("C315").Value = 2
$ws.Range("E315").Value = 0.012
$ws.Range("G315").Value = 0.01

$ws.Range("B316").Value = 2
$ws.Range("C316").Value = 4
$ws.Range("E316").Value = 0.017
$ws.Range("F316").Value = 0.055
$ws.Range("G316").Value = 0.01

$ws.Range("C317").Value = 6
$ws.Range("E317").Value = 0.01
$ws.Range("G317").Value = 0.017

$ws.Range("B318").Value = 2
$ws.Range("C318").Value = 2
$ws.Range("E318").Value = 0.012
$ws.Range("G318").Value = 0.01

$ws.Range("B319").Value = 2
$ws.Range("C319").Value = 4
$ws.Range("E319").Value = 0.017
$ws.Range("F319").Value = 0.061
$ws.Range("G319").Value = 0.01

$ws.Range("C320").Value = 6
$ws.Range("E320").Value = 0.009000000000000001
$ws.Range("G320").Value = 0.017

$ws.Range("B321").Value = 2
$ws.Range("C321").Value = 2
$ws.Range("E321").Value = 0.011000000000000001
$ws.Range("G321").Value = 0.009000000000000001

$ws.Range("B322").Value = 2
$ws.Range("C322").Value = 4
$ws.Range("E322").Value = 0.016
$ws.Range("F322").Value = 0.048
$ws.Range("G322").Value = 0.009000000000000001

$ws.Range("C323").Value = 6
$ws.Range("E323").Value = 0.011000000000000001
$ws.Range("G323").Value = 0.016

$ws.Range("B324").Value = 2
$ws.Range("C324").Value = 2
$ws.Range("E324").Value = 0.013000000000000001
$ws.Range("G324").Value = 0.011000000000000001

$ws.Range("B325").Value = 2
$ws.Range("C325").Value = 4
$ws.Range("E325").Value = 0.018000000000000002
$ws.Range("F325").Value = 0.053
$ws.Range("G325").Value = 0.011000000000000001

# Restore the view state: scroll position / active selection as left by the author
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 282
$ws.Range("C304").Select() | Out-Null
